$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.890.30"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.97%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.514.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.98%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.42%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.97%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("E8").Value = "  -3.24%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.516.50"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.88%  "

# Row 10
$ws.Range("E10").Value = "  -7.22%  "

# Row 12
$ws.Range("E12").Value = "  -4.62%  "

# Row 13
$ws.Range("E13").Value = "  -4.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.977.63"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.723.17"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.81%  "

# Row 16
$ws.Range("E16").Value = "  -4.09%  "

# Row 17
$ws.Range("E17").Value = "  -5.92%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.516.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.76%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.11%  "

# Row 20
$ws.Range("E20").Value = "  -9.21%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.48%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.51%  "

# Row 23
$ws.Range("E23").Value = "  -3.36%  "

# Row 24
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "68.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.10%  "

# Row 26
$ws.Range("E26").Value = "  -6.69%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.29%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.640.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.03%  "

# Row 29
$ws.Range("E29").Value = "  +0.22%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0907"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "478.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.29"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.64%  "

# Row 34
$ws.Range("E34").Value = "  -4.65%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.13%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.76%  "

# Row 37
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.78%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.80%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.68%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.82%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -8.31%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.319"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.49%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.55%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.13%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.17%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.93"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -9.62%  "

# Row 48
$ws.Range("E48").Value = "  -6.81%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.523"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.94%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.56%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.598"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.77%  "
